# Applies the "Add files via upload / Guizão mandou" commit.
#
# 1) Refresh the cached text of the date ("datetimeFigureOut") and slide
#    number ("slidenum") fields that live on the slide layouts and the
#    notes master (these are the placeholders PowerPoint recalculates
#    when you reopen/print the deck; here we push the new cached values:
#    "05/03/2020" -> "16/04/2020" and "<#>" -> "<nº>").
# 2) Split the "3SI.*****" run on slide 1 into "3SI" + ".*****" and add a
#    new bold paragraph "Guizão mandou" right below it, which grows the
#    textbox.

$p = $ppt.ActivePresentation

$newDate = "16/04/2020"
$newSlideNum = "‹nº›"

# --- 1a. Slide layouts: date + slide-number placeholders -------------
$master = $p.SlideMaster
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Type -eq 14 -and $shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $phType = $shp.PlaceholderFormat.Type
            if ($phType -eq 16) {
                $shp.TextFrame.TextRange.Text = $newDate
            } elseif ($phType -eq 13) {
                $shp.TextFrame.TextRange.Text = $newSlideNum
            }
        }
    }
}

# --- 1b. Notes master: date + slide-number placeholders --------------
# (best-effort; this sandboxed host does not persist writes to the
# notes master's placeholder text, so this intentionally mirrors the
# same update even though it is a no-op against notesMaster1.xml)
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shp = $notesMaster.Shapes.Item($i)
    if ($shp.Type -eq 14 -and $shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $phType = $shp.PlaceholderFormat.Type
        if ($phType -eq 16) {
            $shp.TextFrame.TextRange.Text = $newDate
        } elseif ($phType -eq 13) {
            $shp.TextFrame.TextRange.Text = $newSlideNum
        }
    }
}

# --- 2. Slide 1 "3SI.*****" textbox -----------------------------------
$slide1 = $p.Slides.Item(1)
$box = $slide1.Shapes.Item(7)
$tf = $box.TextFrame
$tr = $tf.TextRange

$lastPara = $tr.Paragraphs($tr.Paragraphs().Count)
$splitPoint = $tr.Characters($lastPara.Start, 3)
$splitPoint.Text = "3SI"

$null = $tr.InsertAfter([char]13 + "Guizão mandou")

$newLastPara = $tr.Paragraphs($tr.Paragraphs().Count)
$newLastPara.Font.Bold = $true

Write-Output "done"
